$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 118.36
$ws.Range("I15").Value = 118.36
$ws.Range("K15").Value = 355.08
$ws.Range("M15").Value = -186.08
# Row 19
$ws.Range("H19").Value = 273.07693
$ws.Range("I19").Value = 218.16667
$ws.Range("J19").Value = 320.14285
$ws.Range("K19").Value = 218.16667
$ws.Range("L19").Value = 320.14285
$ws.Range("M19").Value = -43.16667000000001
$ws.Range("N19").Value = -670.14285
# Row 127
$ws.Range("H127").Value = 1379.4
$ws.Range("J127").Value = 1613.4286
$ws.Range("L127").Value = 4840.2858
$ws.Range("N127").Value = -14760.2858
# Row 129
$ws.Range("H129").Value = 244794.69
$ws.Range("J129").Value = 313569.22
$ws.Range("L129").Value = 940707.6599999999
$ws.Range("N129").Value = -950707.6599999999
# Row 135
$ws.Range("I135").Value = 1162.4117
$ws.Range("K135").Value = 10461.7053
$ws.Range("M135").Value = -7926.705300000001
# Row 138
$ws.Range("H138").Value = 3609.93
$ws.Range("I138").Value = 2404.647
$ws.Range("J138").Value = 3856.7952
$ws.Range("K138").Value = 7213.941
$ws.Range("L138").Value = 11570.3856
$ws.Range("M138").Value = -2073.941
$ws.Range("N138").Value = -21850.3856

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 28
$ws.Range("H28").Value = 8145.5713
$ws.Range("I28").Value = 4300
$ws.Range("J28").Value = 17759.5
$ws.Range("K28").Value = 4300
$ws.Range("L28").Value = 17759.5
$ws.Range("M28").Value = -4108
$ws.Range("N28").Value = -18143.5
# Row 32
$ws.Range("H32").Value = 17678.465
$ws.Range("I32").Value = 12239.272
$ws.Range("J32").Value = 35627.8
$ws.Range("K32").Value = 12239.272
$ws.Range("L32").Value = 35627.8
$ws.Range("M32").Value = -11952.272
$ws.Range("N32").Value = -36201.8
# Row 59
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
# Row 61
$ws.Range("H61").Value = 15973372
$ws.Range("I61").Value = 21961880
$ws.Range("J61").Value = 4016.6667
$ws.Range("K61").Value = 21961880
$ws.Range("L61").Value = 4016.6667
$ws.Range("M61").Value = -21961668
$ws.Range("N61").Value = -4440.6667
# Row 74
$ws.Range("H74").Value = 33336580
$ws.Range("I74").Value = 55558572
$ws.Range("K74").Value = 55558572
$ws.Range("M74").Value = -55557698
# Row 77
$ws.Range("H77").Value = 33336580
$ws.Range("I77").Value = 55558572
$ws.Range("K77").Value = 277792860
$ws.Range("M77").Value = -277788492
# Row 99
$ws.Range("H99").Value = 8145.5713
$ws.Range("I99").Value = 4300
$ws.Range("J99").Value = 17759.5
$ws.Range("K99").Value = 4300
$ws.Range("L99").Value = 17759.5
$ws.Range("M99").Value = -1305
$ws.Range("N99").Value = -23749.5
# Row 132
$ws.Range("H132").Value = 7362574
$ws.Range("I132").Value = 8773895
$ws.Range("K132").Value = 26321685
$ws.Range("M132").Value = -26319155
# Row 136
$ws.Range("H136").Value = 15973372
$ws.Range("I136").Value = 21961880
$ws.Range("J136").Value = 4016.6667
$ws.Range("K136").Value = 65885640
$ws.Range("L136").Value = 12050.0001
$ws.Range("M136").Value = -65883090
$ws.Range("N136").Value = -17150.0001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3137.875
$ws.Range("I20").Value = 3238.2307
$ws.Range("K20").Value = 3238.2307
$ws.Range("M20").Value = -2991.2307
# Row 86
$ws.Range("H86").Value = 1956.6
$ws.Range("I86").Value = 1816.2106
$ws.Range("J86").Value = 2401.1667
$ws.Range("K86").Value = 1816.2106
$ws.Range("L86").Value = 2401.1667
$ws.Range("M86").Value = -693.2106000000001
$ws.Range("N86").Value = -4647.1667
# Row 89
$ws.Range("H89").Value = 1956.6
$ws.Range("I89").Value = 1816.2106
$ws.Range("J89").Value = 2401.1667
$ws.Range("K89").Value = 9081.053
$ws.Range("L89").Value = 12005.8335
$ws.Range("M89").Value = -3465.053
$ws.Range("N89").Value = -23237.8335
# Row 134
$ws.Range("H134").Value = 4234.625
$ws.Range("I134").Value = 4018.6562
$ws.Range("J134").Value = 5098.5
$ws.Range("K134").Value = 12055.9686
$ws.Range("L134").Value = 15295.5
$ws.Range("M134").Value = -9520.9686
$ws.Range("N134").Value = -20365.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 50
$ws.Range("H50").Value = 28000
$ws.Range("J50").Value = 28000
$ws.Range("L50").Value = 28000
$ws.Range("N50").Value = -29250
# Row 60
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1259.6383
$ws.Range("I5").Value = 1061.2759
$ws.Range("J5").Value = 1579.2222
$ws.Range("K5").Value = 3183.8277
$ws.Range("L5").Value = 4737.6666
$ws.Range("M5").Value = -3071.8277
$ws.Range("N5").Value = -4961.6666
# Row 34
$ws.Range("H34").Value = 612.8182
$ws.Range("I34").Value = 469.6
$ws.Range("J34").Value = 732.1667
$ws.Range("K34").Value = 1408.8
$ws.Range("L34").Value = 2196.5001
$ws.Range("M34").Value = -1324.8
$ws.Range("N34").Value = -2364.5001
# Row 122
$ws.Range("H122").Value = 1515.0385
$ws.Range("J122").Value = 1587.125
$ws.Range("L122").Value = 14284.125
$ws.Range("N122").Value = -19184.125
# Row 131
$ws.Range("H131").Value = 754.58
$ws.Range("J131").Value = 754.58
$ws.Range("L131").Value = 2263.74
$ws.Range("N131").Value = -12343.74
# Row 135
$ws.Range("H135").Value = 1259.6383
$ws.Range("I135").Value = 1061.2759
$ws.Range("J135").Value = 1579.2222
$ws.Range("K135").Value = 9551.483100000001
$ws.Range("L135").Value = 14212.9998
$ws.Range("M135").Value = -7016.483100000001
$ws.Range("N135").Value = -19282.9998
# Row 137
$ws.Range("H137").Value = 18526160
$ws.Range("J137").Value = 23818490
$ws.Range("L137").Value = 71455470
$ws.Range("N137").Value = -71465670

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3507.111
$ws.Range("I80").Value = 3056.2856
$ws.Range("J80").Value = 3992.6155
$ws.Range("K80").Value = 3056.2856
$ws.Range("L80").Value = 3992.6155
$ws.Range("M80").Value = -2058.2856
$ws.Range("N80").Value = -5988.6155
# Row 83
$ws.Range("H83").Value = 3507.111
$ws.Range("I83").Value = 3056.2856
$ws.Range("J83").Value = 3992.6155
$ws.Range("K83").Value = 15281.428
$ws.Range("L83").Value = 19963.0775
$ws.Range("M83").Value = -10289.428
$ws.Range("N83").Value = -29947.0775
# Row 132
$ws.Range("H132").Value = 4558228
$ws.Range("J132").Value = 49020.547
$ws.Range("L132").Value = 147061.641
$ws.Range("N132").Value = -152121.641

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 550.25
$ws.Range("I16").Value = 580.2857
$ws.Range("K16").Value = 580.2857
$ws.Range("M16").Value = -410.2857
# Row 46
$ws.Range("H46").Value = 2040.4
$ws.Range("J46").Value = 2300.5
$ws.Range("L46").Value = 2300.5
$ws.Range("N46").Value = -2676.5
# Row 93
$ws.Range("H93").Value = 2498.7778
$ws.Range("I93").Value = 2748.1667
$ws.Range("K93").Value = 2748.1667
$ws.Range("M93").Value = -1500.1667
# Row 100
$ws.Range("H100").Value = 2313.7036
$ws.Range("I100").Value = 1833.1818
$ws.Range("J100").Value = 2644.0625
$ws.Range("K100").Value = 1833.1818
$ws.Range("L100").Value = 2644.0625
$ws.Range("M100").Value = -1292.1818
$ws.Range("N100").Value = -3726.0625
# Row 111
$ws.Range("H111").Value = 30387
$ws.Range("J111").Value = 30387
$ws.Range("L111").Value = 30387
$ws.Range("N111").Value = -38567
# Row 115
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
# Row 118
$ws.Range("H118").Value = 30392
$ws.Range("J118").Value = 30392
$ws.Range("L118").Value = 30392
$ws.Range("N118").Value = -33706
